$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_8")

$updates = @{
    8  = 1
    10 = 1
    11 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 4
    24 = 4
    28 = 1
    31 = 1
    38 = 1
    40 = 1
    41 = 2
    47 = 3
    48 = 1
    60 = 1
    62 = 3
    67 = 2
    68 = 2
    69 = 6
    72 = 1
    89 = 1
    90 = 2
    95 = 1
    96 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("Q$row").Value = $updates[$row]
}
